# Weekly update: insert a new record row right after row 13 (before the
# former row 14), shifting all subsequent rows down by one, and populate
# the new row with the latest "Albahaca" price report.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 14; everything currently at row 14
# and below moves down to row 15 and below. Excel carries the formatting
# of the row above (row 13) onto the newly inserted row, which keeps the
# date-formatted style on column D consistent with the rest of the table.
$ws.Rows.Item(14).Insert()

$ws.Cells.Item(14, 1).Value = 8
$ws.Cells.Item(14, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(14, 3).Value = "Coquimbo"
$ws.Cells.Item(14, 4).Value = 44687
$ws.Cells.Item(14, 5).Value = 4
$ws.Cells.Item(14, 6).Value = 100112052
$ws.Cells.Item(14, 7).Value = "Albahaca"
$ws.Cells.Item(14, 8).Value = "Sin especificar"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 1100
$ws.Cells.Item(14, 11).Value = 5000
$ws.Cells.Item(14, 12).Value = 6000
$ws.Cells.Item(14, 13).Value = 5500
$ws.Cells.Item(14, 14).Value = "`$/docena de matas"
$ws.Cells.Item(14, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(14, 16).Value = 917
$ws.Cells.Item(14, 17).Value = 6
$ws.Cells.Item(14, 18).Value = "Hortaliza"
